$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.799.34"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "2.221.81"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.18"
$ws.Range("E5").Value = "  +6.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.87"
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +8.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.58"
$ws.Range("E10").Value = "  +10.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.23"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("E13").Value = "  +7.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.106"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "2.553.37"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.97"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "2.228.23"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "41.689.42"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "0.0₃0964"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.94"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.61"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  +6.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.03"
$ws.Range("E25").Value = "  +9.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  +5.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.80"
$ws.Range("E28").Value = "  +7.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.16"
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  -5.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.78"
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.60"
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.124"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0734"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.90"
$ws.Range("E37").Value = "  +15.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.04"
$ws.Range("E38").Value = "  +9.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0303"
$ws.Range("E39").Value = "  +9.37%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.63"
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.20"
$ws.Range("E43").Value = "  +19.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.203"
$ws.Range("E44").Value = "  +5.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.83"
$ws.Range("E45").Value = "  -2.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.73"
$ws.Range("E46").Value = "  -6.55%  "
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.68"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  +6.01%  "
$ws.Range("E51").Value = "  +0.42%  "
